# Commit: "Mon, Jun 29, 2020  3:06:28 PM"
#
# 1) The table on slide 16 gets a new table style (Design > Table Styles):
#    {4B190747-3347-49FD-9E43-6405A2028CDB} -> {BE39200F-5176-4DA9-8CF9-C34EF5F11CE4}
#
# 2) The deck's theme palette is reset from the "Integral" design back to the
#    stock "Office Theme" colours (table styles are theme-coloured, so the
#    author also restored the default Office palette on the slide master).

function HexToRgbLong([string]$hex) {
    # VBA/COM RGB() packs a colour as R + G*256 + B*65536.
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + $g * 256 + $b * 65536
}

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 16 -------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{BE39200F-5176-4DA9-8CF9-C34EF5F11CE4}")

# --- 2. Restore the default Office Theme colours on the slide master -------
$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$colors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 0; $i -lt $officeColors.Length; $i++) {
    $colors.Item($i + 1).RGB = HexToRgbLong $officeColors[$i]
}
